$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values per "Updated cryptos list on Sun Apr 28 06:08:48 UTC 2024 with GitHub Actions"
# Note: D-column numeric-looking values are prefixed with a literal leading
# apostrophe ('') so Excel stores them as text (matching the original sheet's
# inlineStr/text cells) instead of auto-converting them to numbers.
$ws.Range("D2").Value = '64.012.95'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '3.310.44'
$ws.Range("E3").Value = '  +6.17%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''599.35'
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("D6").Value = '''143.58'
$ws.Range("E6").Value = '  +5.42%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.308.78'
$ws.Range("E8").Value = '  +6.23%  '
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  +2.84%  '
$ws.Range("D11").Value = '''5.48'
$ws.Range("E11").Value = '  +4.70%  '
$ws.Range("E12").Value = '  +3.04%  '
$ws.Range("D13").Value = '''0.0000250'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '''34.92'
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").Value = '3.859.02'
$ws.Range("E15").Value = '  +6.26%  '
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '3.311.74'
$ws.Range("E17").Value = '  +6.14%  '
$ws.Range("D18").Value = '64.067.09'
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("E19").Value = '  +2.73%  '
$ws.Range("D20").Value = '''483.17'
$ws.Range("E20").Value = '  +1.68%  '
$ws.Range("D21").Value = '''14.32'
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("E22").Value = '  +6.19%  '
$ws.Range("E23").Value = '  +3.95%  '
$ws.Range("D24").Value = '''13.60'
$ws.Range("E24").Value = '  +4.60%  '
$ws.Range("D25").Value = '''84.51'
$ws.Range("E25").Value = '  -3.05%  '
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("D28").Value = '''7.28'
$ws.Range("E28").Value = '  +1.66%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '''8.24'
$ws.Range("E30").Value = '  +3.87%  '
$ws.Range("D31").Value = '''2.15'
$ws.Range("E31").Value = '  +4.28%  '
$ws.Range("D32").Value = '''28.67'
$ws.Range("E32").Value = '  +7.15%  '
$ws.Range("E33").Value = '  -0.69%  '
$ws.Range("D34").Value = '''2.56'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").Value = '''1.11'
$ws.Range("E35").Value = '  +2.73%  '
$ws.Range("E36").Value = '  +2.90%  '
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0763'
$ws.Range("E37").Value = '  +6.92%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '''53.31'
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").Value = '''434.03'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("D41").Value = '3.035.05'
$ws.Range("E41").Value = '  +5.18%  '
$ws.Range("E42").Value = '  +4.43%  '
$ws.Range("D43").Value = '''8.45'
$ws.Range("E43").Value = '  +2.22%  '
$ws.Range("E44").Value = '  -5.87%  '
$ws.Range("D45").Value = '''0.269'
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("E46").Value = '  +4.19%  '
$ws.Range("D47").Value = '''26.58'
$ws.Range("E47").Value = '  +3.37%  '
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("E50").Value = '  +2.32%  '
$ws.Range("D51").Value = '''35.30'
$ws.Range("E51").Value = '  +12.08%  '